$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the wording of the existing "Scénarios" journal entry (row 2 / B2)
$ws.Range("B2").Value = "Demande du CdP de refaire les Scénarios car le format n'est pas bon"

# New journal entry in row 3: date + event text, matching the formatting
# already used for the other rows (date style on A, wrapped/centered text on B)
$ws.Range("A3").Value = 43160
$ws.Range("B3").Value = "Le CdP m'a indiqué comment re-faire mes uses cases et scénarios. Priorité mise sur le rendu de ceux-ci"
$ws.Range("B3").HorizontalAlignment = -4108
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("B3").WrapText = $true
$ws.Rows("3:3").RowHeight = 30

# Move the active selection to B6 (next empty row for a future entry)
$null = $ws.Range("B6").Select()

# Restore the saved window x-position captured for this workbook
$excel.ActiveWindow.Left = 2790
